# Add "Group" / ${event.groupName} column (H) to the events export template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the jxls directive comments: lastCell "G7" -> "H7" -----------
$ws.Range("A1").Comment.Text('jx:area(lastCell="H7")') | Out-Null
$ws.Range("A2").Comment.Text('jx:each(items="devices", var="device", lastCell="H7" multisheet="sheetNames")') | Out-Null
$ws.Range("A7").Comment.Text('jx:each(items="device.objects", var="event", lastCell="H7")') | Out-Null

# --- 2. New header cell H6 = "Group" (same look as the other header cells) -
$ws.Range("H6").Value = "Group"
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)

# --- 3. New data cell H7 = "${event.groupName}" (same look as row 7 cells) -
$ws.Range("H7").Value = '${event.groupName}'
$ws.Range("G7").Copy()
$ws.Range("H7").PasteSpecial(-4122)

# --- 4. Column H width, matching the new column added to the report --------
$ws.Columns("H").ColumnWidth = 17.66796875

# --- 5. Restore the cell selection recorded in the saved workbook ----------
$ws.Range("E17").Select() | Out-Null
